$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting rows 42:163 down to 43:164
$ws.Rows("42:42").Insert()

# Fill in the new row 42 with its data. Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T
# carry the same constant values as the rest of the sheet (same market /
# product), while D (fecha), M (volumen), N/O/P (precios) and S (precio $/kg)
# hold the new weekly observation.
$ws.Cells.Item(42, 1).Value = 4
$ws.Cells.Item(42, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(42, 3).Value = "Los Lagos"
$ws.Cells.Item(42, 4).Value = 44497
$ws.Cells.Item(42, 5).Value = 10
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100102
$ws.Cells.Item(42, 8).Value = "Cítricos"
$ws.Cells.Item(42, 9).Value = 100102006
$ws.Cells.Item(42, 10).Value = "Pomelo"
$ws.Cells.Item(42, 11).Value = "Start Ruby"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 13).Value = 160
$ws.Cells.Item(42, 14).Value = 11000
$ws.Cells.Item(42, 15).Value = 12000
$ws.Cells.Item(42, 16).Value = 11500
$ws.Cells.Item(42, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(42, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(42, 19).Value = 821
$ws.Cells.Item(42, 20).Value = 14
